$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new first column ("Sem") ; everything currently in A:J shifts to B:K ---
$ws.Columns("A").Insert()

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Sem"
$ws.Range("A1").HorizontalAlignment = 1   # match the style used by the other header cells
$ws.Range("C1").Value = "CourseName"   # was "courseName" (casing fix)

# --- Row 2 (was row 2 before the shift) ---
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = "intro to fucking"       # was "Intro to Machine Learning X"
$ws.Range("H2").Value = 100                      # was 1100

# --- Row 3 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "EFGH5673"               # was "EFGH5678"
$ws.Range("C3").Value = "intro to puss"          # was "Advanced Introduction to Introduction"

# --- Row 4 ---
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "EFGH5678"               # was "JKEL5679"
$ws.Range("C4").Value = "intro to dick"          # was "Very Unreal Analysis"

# --- Number formatting / alignment for the new "Sem" column ---
$ws.Range("A2:A4").HorizontalAlignment = 1   # xlHAlignGeneral
$ws.Range("A2:A4").NumberFormat = "#,##0"

# --- Row heights grew slightly (18.75 -> 19.5) ---
$ws.Rows("1:4").RowHeight = 19.5

# --- Column widths: the new Sem column and the renamed CourseName column now use the
#     standard column width instead of CourseName's old extra-wide autofit width ---
$ws.Columns("A").ColumnWidth = 11.67
$ws.Columns("C").ColumnWidth = 11.67
